# Mother-In-Law House Expenses - Oct 21 labor + materials update
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Home Summary
# ---------------------------------------------------------------------------
$home = $wb.Worksheets.Item("Home Summary")

$home.Range("B4").Value  = "KES 1,000,000"
$home.Range("B5").Value  = "KES 1,473,573"
$home.Range("B6").Value  = "KES -473,573"
$home.Range("B7").Value  = "'147.36%"
$home.Range("B8").Value  = "KES 16,901"

$home.Range("B12").Value = "KES 15,900"
$home.Range("B13").Value = "KES 107,200"
$home.Range("B14").Value = "KES 1,580,773"
$home.Range("B15").Value = "'158.08%"
$home.Range("B16").Value = "KES -580,773"

$home.Range("B19").Value = "KES 185,705"
$home.Range("B20").Value = "KES 1,766,478"
$home.Range("B21").Value = "KES 766,478"

$home.Range("B25").Value = 485450
$home.Range("C25").Value = 6079.5
$home.Range("D25").Value = 491529.5
$home.Range("E25").Value = "'49.15%"

$home.Range("E26").Value = "'35.26%"
$home.Range("E27").Value = "'28.00%"
$home.Range("E28").Value = "'17.84%"
$home.Range("E29").Value = "'9.15%"
$home.Range("E30").Value = "'3.13%"
$home.Range("E31").Value = "'2.22%"
$home.Range("E32").Value = "'1.08%"
$home.Range("E33").Value = "'1.01%"
$home.Range("E34").Value = "'0.52%"

# ---------------------------------------------------------------------------
# Sheet: Daily Expenses - append Oct 21 entries
# ---------------------------------------------------------------------------
$daily = $wb.Worksheets.Item("Daily Expenses")

# Labor rows (465-468) - clone formatting from the 19/10 unpaid-labor block (455-458)
$daily.Range("A455:I455").Copy($daily.Range("A465:I465"))
$daily.Range("A456:I456").Copy($daily.Range("A466:I466"))
$daily.Range("A457:I457").Copy($daily.Range("A467:I467"))
$daily.Range("A458:I458").Copy($daily.Range("A468:I468"))

$daily.Range("A465").Value = "21/10/2025"
$daily.Range("D465").Value = "Jack - UNPAID"
$daily.Range("E465").Value = 1500

$daily.Range("A466").Value = "21/10/2025"
$daily.Range("D466").Value = "Fundi 1 - UNPAID"
$daily.Range("E466").Value = 1300

$daily.Range("A467").Value = "21/10/2025"
$daily.Range("D467").Value = "Fundi 2 - UNPAID"
$daily.Range("E467").Value = 1300

$daily.Range("A468").Value = "21/10/2025"
$daily.Range("D468").Value = "2 helpers @ 600 each - UNPAID"
$daily.Range("E468").Value = 1200

# Materials rows (469-472) - clone formatting from a paid "Building Materials" row (464)
$daily.Range("A464:I464").Copy($daily.Range("A469:I469"))
$daily.Range("A464:I464").Copy($daily.Range("A470:I470"))
$daily.Range("A464:I464").Copy($daily.Range("A471:I471"))
$daily.Range("A464:I464").Copy($daily.Range("A472:I472"))

$daily.Range("A469").Value = "21/10/2025"
$daily.Range("B469").Value = "Building Materials"
$daily.Range("C469").Value = "Finishing Materials"
$daily.Range("D469").Value = "17ft faceboard 8X1 @ 90"
$daily.Range("E469").Value = 1530
$daily.Range("F469").Value = 25
$daily.Range("G469").Value = 1555
$daily.Range("H469").Value = "Hardware Store"
$daily.Range("I469").Value = "PAID"

$daily.Range("A470").Value = "21/10/2025"
$daily.Range("B470").Value = "Building Materials"
$daily.Range("C470").Value = "Ceiling Materials"
$daily.Range("D470").Value = "6 ceiling board 9mm @ 1050"
$daily.Range("E470").Value = 6300
$daily.Range("F470").Value = 75
$daily.Range("G470").Value = 6375
$daily.Range("H470").Value = "Hardware Store"
$daily.Range("I470").Value = "PAID"

$daily.Range("A471").Value = "21/10/2025"
$daily.Range("B471").Value = "Building Materials"
$daily.Range("C471").Value = "Finishing Materials"
$daily.Range("D471").Value = "8 pcs Aluminium corners strip gold @ 250"
$daily.Range("E471").Value = 2000
$daily.Range("F471").Value = 25
$daily.Range("G471").Value = 2025
$daily.Range("H471").Value = "Hardware Store"
$daily.Range("I471").Value = "PAID"

$daily.Range("A472").Value = "21/10/2025"
$daily.Range("B472").Value = "Building Materials"
$daily.Range("C472").Value = "Ceiling Materials"
$daily.Range("D472").Value = "Spacers 2.5mm 4 packets @ 100"
$daily.Range("E472").Value = 400
$daily.Range("F472").Value = 5
$daily.Range("G472").Value = 405
$daily.Range("H472").Value = "Hardware Store"
$daily.Range("I472").Value = "PAID"

# ---------------------------------------------------------------------------
# Sheet: M-Pesa Fees
# ---------------------------------------------------------------------------
$mpesa = $wb.Worksheets.Item("M-Pesa Fees")

$mpesa.Range("C6").Value  = 66
$mpesa.Range("D6").Value  = 1650
$mpesa.Range("C8").Value  = 18
$mpesa.Range("D8").Value  = 1350
$mpesa.Range("C11").Value = 156
$mpesa.Range("C15").Value = 57
$mpesa.Range("D15").Value = 285
$mpesa.Range("B20").Value = "KES 16,901"

# ---------------------------------------------------------------------------
# Sheet: Unpaid Labor - append Oct 21 unpaid labor and move totals row
# ---------------------------------------------------------------------------
$unpaid = $wb.Worksheets.Item("Unpaid Labor")

# Move the "Total Unpaid Labor" row from row 13 down to row 17 (clone formatting first)
$unpaid.Range("A13:D13").Copy($unpaid.Range("A17:D17"))
$unpaid.Range("B17").Value = ""
$unpaid.Range("C17").Value = "KES 15,900"

# New data rows 12-15, cloned from an existing "PENDING" row (row 11)
$unpaid.Range("A11:D11").Copy($unpaid.Range("A12:D12"))
$unpaid.Range("A11:D11").Copy($unpaid.Range("A13:D13"))
$unpaid.Range("A11:D11").Copy($unpaid.Range("A14:D14"))
$unpaid.Range("A11:D11").Copy($unpaid.Range("A15:D15"))

$unpaid.Range("A12").Value = "21/10/2025"
$unpaid.Range("B12").Value = "Jack - UNPAID"
$unpaid.Range("C12").Value = 1500

$unpaid.Range("A13").Value = "21/10/2025"
$unpaid.Range("B13").Value = "Fundi 1 - UNPAID"
$unpaid.Range("C13").Value = 1300

$unpaid.Range("A14").Value = "21/10/2025"
$unpaid.Range("B14").Value = "Fundi 2 - UNPAID"
$unpaid.Range("C14").Value = 1300

$unpaid.Range("A15").Value = "21/10/2025"
$unpaid.Range("B15").Value = "2 helpers @ 600 each - UNPAID"
$unpaid.Range("C15").Value = 1200

# ---------------------------------------------------------------------------
# Sheet: Pending Purchases
# ---------------------------------------------------------------------------
$pending = $wb.Worksheets.Item("Pending Purchases")

$pending.Range("C15").Value = 38205
$pending.Range("C17").Value = "KES 185,705"
$pending.Range("B20").Value = "KES 1,473,573"
$pending.Range("B22").Value = "KES 15,900"
$pending.Range("B23").Value = "KES 185,705"
$pending.Range("B25").Value = "KES 1,766,478"
$pending.Range("B26").Value = "KES 1,000,000"
$pending.Range("B27").Value = "KES 766,478"

Write-Host "Edit complete"
